$d = $word.ActiveDocument

# 1) "Klasa abstrakcyjna opisująca poszczególnych użytkowników systemu." ->
#    "Klasa opisująca poszczególnych użytkowników systemu."
$d.Content.Find.Execute(
    "Klasa abstrakcyjna opisująca poszczególnych użytkowników systemu.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Klasa opisująca poszczególnych użytkowników systemu.", 2)

# 2) "Klasa abstrakcyjna opisująca pracowników firmy. Dziedziczy z klasy " ->
#    "Klasa opisująca pracowników firmy. Dziedziczy z klasy "
$d.Content.Find.Execute(
    "Klasa abstrakcyjna opisująca pracowników firmy. Dziedziczy z klasy ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Klasa opisująca pracowników firmy. Dziedziczy z klasy ", 2)

# 3) Remove the "Ogólny opis" heading paragraph entirely.
$rngOgolny = $d.Content
$rngOgolny.Find.Execute("Ogólny opis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngOgolny.Paragraphs(1).Range.Delete()

# 4) Remove the "Opis pól i metod" heading paragraph entirely.
$rngOpis = $d.Content
$rngOpis.Find.Execute("Opis pól i metod", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngOpis.Paragraphs(1).Range.Delete()

# 5) Merge the "na " run with the "wyszukanie w systemie informacji o danym
#    pracowniku." run (which carries a lastRenderedPageBreak) into a single
#    run, removing the page-break marker, while leaving the preceding
#    "właścicielowi " run untouched. Locating the match narrows the range to
#    exactly that run's text; re-running Find/Replace on that same narrowed
#    range (with identical replacement text) forces the engine to rebuild the
#    run, absorbing the adjoining "na " text and dropping the page break,
#    without disturbing the separate "właścicielowi " run before it.
$rngWysz = $d.Content
$rngWysz.Find.Execute(
    "wyszukanie w systemie informacji o danym pracowniku.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngWysz.Find.Execute(
    "wyszukanie w systemie informacji o danym pracowniku.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "wyszukanie w systemie informacji o danym pracowniku.", 2)

Write-Host "All edits applied"
